$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Spon2"
$ws.Cells.Item(2, 3).Value = "Itga5"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.3355466666666667
$ws.Cells.Item(2, 8).Value = 1.00664
$ws.Cells.Item(2, 9).Value = 0.02388108083384315
$ws.Cells.Item(2, 10).Value = 0.02388108083384315
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 34.07074633333333
$ws.Cells.Item(2, 14).Value = 102.212239
$ws.Cells.Item(2, 15).Value = 0.5171464495142372
$ws.Cells.Item(2, 16).Value = 0.5171464495142373
$ws.Cells.Item(2, 17).Value = 11.43232536299555
$ws.Cells.Item(2, 18).Value = 102.89092826696
$ws.Cells.Item(2, 19).Value = 0.01235001616378448
$ws.Cells.Item(2, 20).Value = 0.01235001616378449

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Spon2"
$ws.Cells.Item(3, 3).Value = "Itga5"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.3355466666666667
$ws.Cells.Item(3, 8).Value = 1.00664
$ws.Cells.Item(3, 9).Value = 0.02388108083384315
$ws.Cells.Item(3, 10).Value = 0.02388108083384315
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 27.685497
$ws.Cells.Item(3, 14).Value = 83.056491
$ws.Cells.Item(3, 15).Value = 0.420227262899125
$ws.Cells.Item(3, 16).Value = 0.4202272628991251
$ws.Cells.Item(3, 17).Value = 9.28977623336
$ws.Cells.Item(3, 18).Value = 83.60798610024
$ws.Cells.Item(3, 19).Value = 0.01003548123387866
$ws.Cells.Item(3, 20).Value = 0.01003548123387866

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Spon2"
$ws.Cells.Item(4, 3).Value = "Itga5"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.3355466666666667
$ws.Cells.Item(4, 8).Value = 1.00664
$ws.Cells.Item(4, 9).Value = 0.02388108083384315
$ws.Cells.Item(4, 10).Value = 0.02388108083384315
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.125957666666666
$ws.Cells.Item(4, 14).Value = 12.377873
$ws.Cells.Item(4, 15).Value = 0.06262628758663766
$ws.Cells.Item(4, 16).Value = 0.06262628758663766
$ws.Cells.Item(4, 17).Value = 1.384451341857778
$ws.Cells.Item(4, 18).Value = 12.46006207672
$ws.Cells.Item(4, 19).Value = 0.001495583436180002
$ws.Cells.Item(4, 20).Value = 0.001495583436180002

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Spon2"
$ws.Cells.Item(5, 3).Value = "Itga5"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 12.878362
$ws.Cells.Item(5, 8).Value = 38.635086
$ws.Cells.Item(5, 9).Value = 0.9165616424823987
$ws.Cells.Item(5, 10).Value = 0.9165616424823989
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 34.07074633333333
$ws.Cells.Item(5, 14).Value = 102.212239
$ws.Cells.Item(5, 15).Value = 0.5171464495142372
$ws.Cells.Item(5, 16).Value = 0.5171464495142373
$ws.Cells.Item(5, 17).Value = 438.7754048908392
$ws.Cells.Item(5, 18).Value = 3948.978644017554
$ws.Cells.Item(5, 19).Value = 0.4739965991707101
$ws.Cells.Item(5, 20).Value = 0.4739965991707103

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Spon2"
$ws.Cells.Item(6, 3).Value = "Itga5"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 12.878362
$ws.Cells.Item(6, 8).Value = 38.635086
$ws.Cells.Item(6, 9).Value = 0.9165616424823987
$ws.Cells.Item(6, 10).Value = 0.9165616424823989
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 27.685497
$ws.Cells.Item(6, 14).Value = 83.056491
$ws.Cells.Item(6, 15).Value = 0.420227262899125
$ws.Cells.Item(6, 16).Value = 0.4202272628991251
$ws.Cells.Item(6, 17).Value = 356.543852515914
$ws.Cells.Item(6, 18).Value = 3208.894672643226
$ws.Cells.Item(6, 19).Value = 0.3851641902987048
$ws.Cells.Item(6, 20).Value = 0.385164190298705

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Spon2"
$ws.Cells.Item(7, 3).Value = "Itga5"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 12.878362
$ws.Cells.Item(7, 8).Value = 38.635086
$ws.Cells.Item(7, 9).Value = 0.9165616424823987
$ws.Cells.Item(7, 10).Value = 0.9165616424823989
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 4.125957666666666
$ws.Cells.Item(7, 14).Value = 12.377873
$ws.Cells.Item(7, 15).Value = 0.06262628758663766
$ws.Cells.Item(7, 16).Value = 0.06262628758663766
$ws.Cells.Item(7, 17).Value = 53.13557642800867
$ws.Cells.Item(7, 18).Value = 478.220187852078
$ws.Cells.Item(7, 19).Value = 0.05740085301298366
$ws.Cells.Item(7, 20).Value = 0.05740085301298368

$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Spon2"
$ws.Cells.Item(8, 3).Value = "Itga5"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.8368233333333333
$ws.Cells.Item(8, 8).Value = 2.51047
$ws.Cells.Item(8, 9).Value = 0.05955727668375805
$ws.Cells.Item(8, 10).Value = 0.05955727668375806
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 34.07074633333333
$ws.Cells.Item(8, 14).Value = 102.212239
$ws.Cells.Item(8, 15).Value = 0.5171464495142372
$ws.Cells.Item(8, 16).Value = 0.5171464495142373
$ws.Cells.Item(8, 17).Value = 28.51119551581444
$ws.Cells.Item(8, 18).Value = 256.6007596423299
$ws.Cells.Item(8, 19).Value = 0.03079983417974253
$ws.Cells.Item(8, 20).Value = 0.03079983417974255

$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Spon2"
$ws.Cells.Item(9, 3).Value = "Itga5"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.8368233333333333
$ws.Cells.Item(9, 8).Value = 2.51047
$ws.Cells.Item(9, 9).Value = 0.05955727668375805
$ws.Cells.Item(9, 10).Value = 0.05955727668375806
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 27.685497
$ws.Cells.Item(9, 14).Value = 83.056491
$ws.Cells.Item(9, 15).Value = 0.420227262899125
$ws.Cells.Item(9, 16).Value = 0.4202272628991251
$ws.Cells.Item(9, 17).Value = 23.16786988453
$ws.Cells.Item(9, 18).Value = 208.51082896077
$ws.Cells.Item(9, 19).Value = 0.02502759136654152
$ws.Cells.Item(9, 20).Value = 0.02502759136654153

$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Spon2"
$ws.Cells.Item(10, 3).Value = "Itga5"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.8368233333333333
$ws.Cells.Item(10, 8).Value = 2.51047
$ws.Cells.Item(10, 9).Value = 0.05955727668375805
$ws.Cells.Item(10, 10).Value = 0.05955727668375806
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 4.125957666666666
$ws.Cells.Item(10, 14).Value = 12.377873
$ws.Cells.Item(10, 15).Value = 0.06262628758663766
$ws.Cells.Item(10, 16).Value = 0.06262628758663766
$ws.Cells.Item(10, 17).Value = 3.452697647812222
$ws.Cells.Item(10, 18).Value = 31.07427883031
$ws.Cells.Item(10, 19).Value = 0.003729851137473981
$ws.Cells.Item(10, 20).Value = 0.003729851137473982
